$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: Insert two new columns at D:E, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert()

# Step 2: Copy formats from the (now-shifted) F:M columns into new D:E columns
# so the new columns inherit the same number formats/styles as the rest of the table.
$ws.Range("F8:M102").Copy()
$ws.Range("D8:E102").PasteSpecial(-4122)
$ws.Range("F7:M7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F38:M38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80:M80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Write the restated financial data (10 quarterly columns D..M) for every data row.
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643
$ws.Cells.Item(8, 4).Value = 2229300
$ws.Cells.Item(8, 5).Value = 2047000
$ws.Cells.Item(8, 6).Value = 1953900
$ws.Cells.Item(8, 7).Value = 1900400
$ws.Cells.Item(8, 8).Value = 2012600
$ws.Cells.Item(8, 9).Value = 1886700
$ws.Cells.Item(8, 10).Value = 1895900
$ws.Cells.Item(8, 11).Value = 1891700
$ws.Cells.Item(8, 12).Value = 1950000
$ws.Cells.Item(8, 13).Value = 1923200
$ws.Cells.Item(9, 4).Value = 1893200
$ws.Cells.Item(9, 5).Value = 1737700
$ws.Cells.Item(9, 6).Value = 1663000
$ws.Cells.Item(9, 7).Value = 1631300
$ws.Cells.Item(9, 8).Value = 1701500
$ws.Cells.Item(9, 9).Value = 1591600
$ws.Cells.Item(9, 10).Value = 1621400
$ws.Cells.Item(9, 11).Value = 1625400
$ws.Cells.Item(9, 12).Value = 1678000
$ws.Cells.Item(9, 13).Value = 1655100
$ws.Cells.Item(10, 4).Value = 336100
$ws.Cells.Item(10, 5).Value = 309300
$ws.Cells.Item(10, 6).Value = 290900
$ws.Cells.Item(10, 7).Value = 269100
$ws.Cells.Item(10, 8).Value = 311100
$ws.Cells.Item(10, 9).Value = 295100
$ws.Cells.Item(10, 10).Value = 274500
$ws.Cells.Item(10, 11).Value = 266300
$ws.Cells.Item(10, 12).Value = 272000
$ws.Cells.Item(10, 13).Value = 268100
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(14, 4).Value = 1600
$ws.Cells.Item(14, 5).Value = 200
$ws.Cells.Item(14, 6).Value = 1300
$ws.Cells.Item(14, 7).Value = 100
$ws.Cells.Item(14, 8).Value = 58400
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 600
$ws.Cells.Item(14, 12).Value = 2600
$ws.Cells.Item(14, 13).Value = 500
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(17, 4).Value = 2115700
$ws.Cells.Item(17, 5).Value = 1935300
$ws.Cells.Item(17, 6).Value = 1854200
$ws.Cells.Item(17, 7).Value = 1822400
$ws.Cells.Item(17, 8).Value = 1964600
$ws.Cells.Item(17, 9).Value = 1780600
$ws.Cells.Item(17, 10).Value = 1803500
$ws.Cells.Item(17, 11).Value = 1809000
$ws.Cells.Item(17, 12).Value = 1875500
$ws.Cells.Item(17, 13).Value = 1837100
$ws.Cells.Item(18, 4).Value = 113600
$ws.Cells.Item(18, 5).Value = 111700
$ws.Cells.Item(18, 6).Value = 99700
$ws.Cells.Item(18, 7).Value = 78000
$ws.Cells.Item(18, 8).Value = 48000
$ws.Cells.Item(18, 9).Value = 106100
$ws.Cells.Item(18, 10).Value = 92400
$ws.Cells.Item(18, 11).Value = 82700
$ws.Cells.Item(18, 12).Value = 74500
$ws.Cells.Item(18, 13).Value = 86100
$ws.Cells.Item(20, 4).Value = 1400
$ws.Cells.Item(20, 5).Value = 1500
$ws.Cells.Item(20, 6).Value = 1300
$ws.Cells.Item(20, 7).Value = 1300
$ws.Cells.Item(20, 8).Value = 800
$ws.Cells.Item(20, 9).Value = 700
$ws.Cells.Item(20, 10).Value = 500
$ws.Cells.Item(20, 11).Value = 300
$ws.Cells.Item(20, 12).Value = 100
$ws.Cells.Item(20, 13).Value = 100
$ws.Cells.Item(21, 4).Value = 135800
$ws.Cells.Item(21, 5).Value = 132800
$ws.Cells.Item(21, 6).Value = 121200
$ws.Cells.Item(21, 7).Value = 99700
$ws.Cells.Item(21, 8).Value = 70800
$ws.Cells.Item(21, 9).Value = 128700
$ws.Cells.Item(21, 10).Value = 115100
$ws.Cells.Item(21, 11).Value = 105400
$ws.Cells.Item(21, 12).Value = 94600
$ws.Cells.Item(21, 13).Value = 106700
$ws.Cells.Item(22, 4).Value = 3500
$ws.Cells.Item(22, 5).Value = 3600
$ws.Cells.Item(22, 6).Value = 3500
$ws.Cells.Item(22, 7).Value = 3000
$ws.Cells.Item(22, 8).Value = 3300
$ws.Cells.Item(22, 9).Value = 3300
$ws.Cells.Item(22, 10).Value = 3100
$ws.Cells.Item(22, 11).Value = 3100
$ws.Cells.Item(22, 12).Value = 3700
$ws.Cells.Item(22, 13).Value = 3500
$ws.Cells.Item(23, 4).Value = 111500
$ws.Cells.Item(23, 5).Value = 109700
$ws.Cells.Item(23, 6).Value = 97600
$ws.Cells.Item(23, 7).Value = 76300
$ws.Cells.Item(23, 8).Value = 45600
$ws.Cells.Item(23, 9).Value = 103400
$ws.Cells.Item(23, 10).Value = 89800
$ws.Cells.Item(23, 11).Value = 80000
$ws.Cells.Item(23, 12).Value = 71000
$ws.Cells.Item(23, 13).Value = 82700
$ws.Cells.Item(24, 4).Value = 32300
$ws.Cells.Item(24, 5).Value = 29600
$ws.Cells.Item(24, 6).Value = 26500
$ws.Cells.Item(24, 7).Value = 20600
$ws.Cells.Item(24, 8).Value = 31600
$ws.Cells.Item(24, 9).Value = 38600
$ws.Cells.Item(24, 10).Value = 33000
$ws.Cells.Item(24, 11).Value = 26800
$ws.Cells.Item(24, 12).Value = 28500
$ws.Cells.Item(24, 13).Value = 30800
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 4).Value = 79200
$ws.Cells.Item(26, 5).Value = 80000
$ws.Cells.Item(26, 6).Value = 71000
$ws.Cells.Item(26, 7).Value = 55700
$ws.Cells.Item(26, 8).Value = 14000
$ws.Cells.Item(26, 9).Value = 64800
$ws.Cells.Item(26, 10).Value = 56800
$ws.Cells.Item(26, 11).Value = 53100
$ws.Cells.Item(26, 12).Value = 42400
$ws.Cells.Item(26, 13).Value = 52000
$ws.Cells.Item(27, 4).Value = 79200
$ws.Cells.Item(27, 5).Value = 80000
$ws.Cells.Item(27, 6).Value = 71000
$ws.Cells.Item(27, 7).Value = 55700
$ws.Cells.Item(27, 8).Value = 14000
$ws.Cells.Item(27, 9).Value = 64800
$ws.Cells.Item(27, 10).Value = 56800
$ws.Cells.Item(27, 11).Value = 53100
$ws.Cells.Item(27, 12).Value = 42200
$ws.Cells.Item(27, 13).Value = 51900
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(29, 4).Value = -1200
$ws.Cells.Item(29, 5).Value = -600
$ws.Cells.Item(29, 6).Value = -200
$ws.Cells.Item(29, 7).Value = -300
$ws.Cells.Item(29, 8).Value = 39200
$ws.Cells.Item(29, 9).Value = -200
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = -500
$ws.Cells.Item(29, 12).Value = -1600
$ws.Cells.Item(29, 13).Value = -400
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(32, 4).Value = -1400
$ws.Cells.Item(32, 5).Value = -1500
$ws.Cells.Item(32, 6).Value = -1300
$ws.Cells.Item(32, 7).Value = -1300
$ws.Cells.Item(32, 8).Value = -800
$ws.Cells.Item(32, 9).Value = -700
$ws.Cells.Item(32, 10).Value = -500
$ws.Cells.Item(32, 11).Value = -300
$ws.Cells.Item(32, 12).Value = -100
$ws.Cells.Item(32, 13).Value = -100
$ws.Cells.Item(33, 4).Value = 78000
$ws.Cells.Item(33, 5).Value = 79400
$ws.Cells.Item(33, 6).Value = 70800
$ws.Cells.Item(33, 7).Value = 55400
$ws.Cells.Item(33, 8).Value = 53200
$ws.Cells.Item(33, 9).Value = 64600
$ws.Cells.Item(33, 10).Value = 56800
$ws.Cells.Item(33, 11).Value = 52600
$ws.Cells.Item(33, 12).Value = 40700
$ws.Cells.Item(33, 13).Value = 51500
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(35, 4).Value = 78000
$ws.Cells.Item(35, 5).Value = 79400
$ws.Cells.Item(35, 6).Value = 70800
$ws.Cells.Item(35, 7).Value = 55400
$ws.Cells.Item(35, 8).Value = 53200
$ws.Cells.Item(35, 9).Value = 64600
$ws.Cells.Item(35, 10).Value = 56800
$ws.Cells.Item(35, 11).Value = 52600
$ws.Cells.Item(35, 12).Value = 40700
$ws.Cells.Item(35, 13).Value = 51500
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643
$ws.Cells.Item(41, 4).Value = 363900
$ws.Cells.Item(41, 5).Value = 353900
$ws.Cells.Item(41, 6).Value = 306600
$ws.Cells.Item(41, 7).Value = 352400
$ws.Cells.Item(41, 8).Value = 467400
$ws.Cells.Item(41, 9).Value = 480500
$ws.Cells.Item(41, 10).Value = 385400
$ws.Cells.Item(41, 11).Value = 302800
$ws.Cells.Item(41, 12).Value = 464600
$ws.Cells.Item(41, 13).Value = 504600
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 0
$ws.Cells.Item(43, 4).Value = 1931900
$ws.Cells.Item(43, 5).Value = 1884100
$ws.Cells.Item(43, 6).Value = 1791400
$ws.Cells.Item(43, 7).Value = 1731800
$ws.Cells.Item(43, 8).Value = 1730500
$ws.Cells.Item(43, 9).Value = 1673400
$ws.Cells.Item(43, 10).Value = 1633300
$ws.Cells.Item(43, 11).Value = 1508500
$ws.Cells.Item(43, 12).Value = 1495400
$ws.Cells.Item(43, 13).Value = 1495600
$ws.Cells.Item(44, 4).Value = 42300
$ws.Cells.Item(44, 5).Value = 41900
$ws.Cells.Item(44, 6).Value = 43700
$ws.Cells.Item(44, 7).Value = 43700
$ws.Cells.Item(44, 8).Value = 42700
$ws.Cells.Item(44, 9).Value = 43700
$ws.Cells.Item(44, 10).Value = 42500
$ws.Cells.Item(44, 11).Value = 187600
$ws.Cells.Item(44, 12).Value = 168100
$ws.Cells.Item(44, 13).Value = 186500
$ws.Cells.Item(45, 4).Value = 48100
$ws.Cells.Item(45, 5).Value = 44200
$ws.Cells.Item(45, 6).Value = 48500
$ws.Cells.Item(45, 7).Value = 41800
$ws.Cells.Item(45, 8).Value = 43800
$ws.Cells.Item(45, 9).Value = 34400
$ws.Cells.Item(45, 10).Value = 38700
$ws.Cells.Item(45, 11).Value = 43000
$ws.Cells.Item(45, 12).Value = 40900
$ws.Cells.Item(45, 13).Value = 60500
$ws.Cells.Item(46, 4).Value = 2386200
$ws.Cells.Item(46, 5).Value = 2324200
$ws.Cells.Item(46, 6).Value = 2190200
$ws.Cells.Item(46, 7).Value = 2169800
$ws.Cells.Item(46, 8).Value = 2284500
$ws.Cells.Item(46, 9).Value = 2232000
$ws.Cells.Item(46, 10).Value = 2099800
$ws.Cells.Item(46, 11).Value = 2041800
$ws.Cells.Item(46, 12).Value = 2169100
$ws.Cells.Item(46, 13).Value = 2247100
$ws.Cells.Item(47, 4).Value = 2900
$ws.Cells.Item(47, 5).Value = 3500
$ws.Cells.Item(47, 6).Value = 4200
$ws.Cells.Item(47, 7).Value = 3700
$ws.Cells.Item(47, 8).Value = 2300
$ws.Cells.Item(47, 9).Value = 7600
$ws.Cells.Item(47, 10).Value = 7700
$ws.Cells.Item(47, 11).Value = 9200
$ws.Cells.Item(47, 12).Value = 8800
$ws.Cells.Item(47, 13).Value = 8200
$ws.Cells.Item(48, 4).Value = 134400
$ws.Cells.Item(48, 5).Value = 128800
$ws.Cells.Item(48, 6).Value = 125000
$ws.Cells.Item(48, 7).Value = 124800
$ws.Cells.Item(48, 8).Value = 127200
$ws.Cells.Item(48, 9).Value = 126600
$ws.Cells.Item(48, 10).Value = 127800
$ws.Cells.Item(48, 11).Value = 131500
$ws.Cells.Item(48, 12).Value = 128000
$ws.Cells.Item(48, 13).Value = 128300
$ws.Cells.Item(49, 4).Value = 1479200
$ws.Cells.Item(49, 5).Value = 1449700
$ws.Cells.Item(49, 6).Value = 1459900
$ws.Cells.Item(49, 7).Value = 1449400
$ws.Cells.Item(49, 8).Value = 1459900
$ws.Cells.Item(49, 9).Value = 1507700
$ws.Cells.Item(49, 10).Value = 1519700
$ws.Cells.Item(49, 11).Value = 1531800
$ws.Cells.Item(49, 12).Value = 1467000
$ws.Cells.Item(49, 13).Value = 1479400
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(52, 4).Value = 86200
$ws.Cells.Item(52, 5).Value = 89700
$ws.Cells.Item(52, 6).Value = 90100
$ws.Cells.Item(52, 7).Value = 94300
$ws.Cells.Item(52, 8).Value = 92000
$ws.Cells.Item(52, 9).Value = 92500
$ws.Cells.Item(52, 10).Value = 92500
$ws.Cells.Item(52, 11).Value = 90900
$ws.Cells.Item(52, 12).Value = 79600
$ws.Cells.Item(52, 13).Value = 31900
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(54, 4).Value = 4088800
$ws.Cells.Item(54, 5).Value = 3995900
$ws.Cells.Item(54, 6).Value = 3869400
$ws.Cells.Item(54, 7).Value = 3842000
$ws.Cells.Item(54, 8).Value = 3965900
$ws.Cells.Item(54, 9).Value = 3966300
$ws.Cells.Item(54, 10).Value = 3847600
$ws.Cells.Item(54, 11).Value = 3805200
$ws.Cells.Item(54, 12).Value = 3852400
$ws.Cells.Item(54, 13).Value = 3895000
$ws.Cells.Item(57, 4).Value = 652100
$ws.Cells.Item(57, 5).Value = 529600
$ws.Cells.Item(57, 6).Value = 505400
$ws.Cells.Item(57, 7).Value = 487300
$ws.Cells.Item(57, 8).Value = 567800
$ws.Cells.Item(57, 9).Value = 485700
$ws.Cells.Item(57, 10).Value = 469300
$ws.Cells.Item(57, 11).Value = 467000
$ws.Cells.Item(57, 12).Value = 501200
$ws.Cells.Item(57, 13).Value = 465400
$ws.Cells.Item(58, 4).Value = 16000
$ws.Cells.Item(58, 5).Value = 15800
$ws.Cells.Item(58, 6).Value = 15600
$ws.Cells.Item(58, 7).Value = 15500
$ws.Cells.Item(58, 8).Value = 15400
$ws.Cells.Item(58, 9).Value = 15400
$ws.Cells.Item(58, 10).Value = 15300
$ws.Cells.Item(58, 11).Value = 15300
$ws.Cells.Item(58, 12).Value = 15000
$ws.Cells.Item(58, 13).Value = 19900
$ws.Cells.Item(59, 4).Value = 1066300
$ws.Cells.Item(59, 5).Value = 1044200
$ws.Cells.Item(59, 6).Value = 983900
$ws.Cells.Item(59, 7).Value = 1009500
$ws.Cells.Item(59, 8).Value = 1067700
$ws.Cells.Item(59, 9).Value = 1051200
$ws.Cells.Item(59, 10).Value = 984200
$ws.Cells.Item(59, 11).Value = 990600
$ws.Cells.Item(59, 12).Value = 995500
$ws.Cells.Item(59, 13).Value = 961600
$ws.Cells.Item(60, 4).Value = 1734400
$ws.Cells.Item(60, 5).Value = 1589600
$ws.Cells.Item(60, 6).Value = 1504900
$ws.Cells.Item(60, 7).Value = 1512300
$ws.Cells.Item(60, 8).Value = 1651000
$ws.Cells.Item(60, 9).Value = 1552300
$ws.Cells.Item(60, 10).Value = 1468900
$ws.Cells.Item(60, 11).Value = 1472900
$ws.Cells.Item(60, 12).Value = 1511800
$ws.Cells.Item(60, 13).Value = 1446800
$ws.Cells.Item(61, 4).Value = 279800
$ws.Cells.Item(61, 5).Value = 283600
$ws.Cells.Item(61, 6).Value = 287500
$ws.Cells.Item(61, 7).Value = 291500
$ws.Cells.Item(61, 8).Value = 294800
$ws.Cells.Item(61, 9).Value = 398500
$ws.Cells.Item(61, 10).Value = 401900
$ws.Cells.Item(61, 11).Value = 405400
$ws.Cells.Item(61, 12).Value = 408300
$ws.Cells.Item(61, 13).Value = 503400
$ws.Cells.Item(62, 4).Value = 333200
$ws.Cells.Item(62, 5).Value = 343600
$ws.Cells.Item(62, 6).Value = 341800
$ws.Cells.Item(62, 7).Value = 346600
$ws.Cells.Item(62, 8).Value = 346000
$ws.Cells.Item(62, 9).Value = 397700
$ws.Cells.Item(62, 10).Value = 397900
$ws.Cells.Item(62, 11).Value = 392800
$ws.Cells.Item(62, 12).Value = 394400
$ws.Cells.Item(62, 13).Value = 360800
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(66, 4).Value = 2348300
$ws.Cells.Item(66, 5).Value = 2217700
$ws.Cells.Item(66, 6).Value = 2135100
$ws.Cells.Item(66, 7).Value = 2151200
$ws.Cells.Item(66, 8).Value = 2292600
$ws.Cells.Item(66, 9).Value = 2349300
$ws.Cells.Item(66, 10).Value = 2269500
$ws.Cells.Item(66, 11).Value = 2271800
$ws.Cells.Item(66, 12).Value = 2315300
$ws.Cells.Item(66, 13).Value = 2312400
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(72, 4).Value = 2060400
$ws.Cells.Item(72, 5).Value = 1987100
$ws.Cells.Item(72, 6).Value = 1912400
$ws.Cells.Item(72, 7).Value = 1846300
$ws.Cells.Item(72, 8).Value = 1796600
$ws.Cells.Item(72, 9).Value = 1748100
$ws.Cells.Item(72, 10).Value = 1688300
$ws.Cells.Item(72, 11).Value = 1644100
$ws.Cells.Item(72, 12).Value = 1596300
$ws.Cells.Item(72, 13).Value = 1560500
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0
$ws.Cells.Item(76, 4).Value = 1740500
$ws.Cells.Item(76, 5).Value = 1778200
$ws.Cells.Item(76, 6).Value = 1734300
$ws.Cells.Item(76, 7).Value = 1690800
$ws.Cells.Item(76, 8).Value = 1673300
$ws.Cells.Item(76, 9).Value = 1617000
$ws.Cells.Item(76, 10).Value = 1578100
$ws.Cells.Item(76, 11).Value = 1533300
$ws.Cells.Item(76, 12).Value = 1537100
$ws.Cells.Item(76, 13).Value = 1582600
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643
$ws.Cells.Item(81, 4).Value = 78000
$ws.Cells.Item(81, 5).Value = 79400
$ws.Cells.Item(81, 6).Value = 70800
$ws.Cells.Item(81, 7).Value = 55400
$ws.Cells.Item(81, 8).Value = 53200
$ws.Cells.Item(81, 9).Value = 64600
$ws.Cells.Item(81, 10).Value = 56800
$ws.Cells.Item(81, 11).Value = 52600
$ws.Cells.Item(81, 12).Value = 40700
$ws.Cells.Item(81, 13).Value = 51500
$ws.Cells.Item(83, 4).Value = 20700
$ws.Cells.Item(83, 5).Value = 19600
$ws.Cells.Item(83, 6).Value = 20200
$ws.Cells.Item(83, 7).Value = 20400
$ws.Cells.Item(83, 8).Value = 22000
$ws.Cells.Item(83, 9).Value = 21900
$ws.Cells.Item(83, 10).Value = 22200
$ws.Cells.Item(83, 11).Value = 22400
$ws.Cells.Item(83, 12).Value = 20000
$ws.Cells.Item(83, 13).Value = 20500
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(89, 4).Value = 205100
$ws.Cells.Item(89, 5).Value = 98500
$ws.Cells.Item(89, 6).Value = 26400
$ws.Cells.Item(89, 7).Value = -59100
$ws.Cells.Item(89, 8).Value = 127800
$ws.Cells.Item(89, 9).Value = 135400
$ws.Cells.Item(89, 10).Value = 106800
$ws.Cells.Item(89, 11).Value = -3900
$ws.Cells.Item(89, 12).Value = 135600
$ws.Cells.Item(89, 13).Value = 81100
$ws.Cells.Item(91, 4).Value = -14800
$ws.Cells.Item(91, 5).Value = -12800
$ws.Cells.Item(91, 6).Value = -9300
$ws.Cells.Item(91, 7).Value = -6600
$ws.Cells.Item(91, 8).Value = -8600
$ws.Cells.Item(91, 9).Value = -8400
$ws.Cells.Item(91, 10).Value = -7100
$ws.Cells.Item(91, 11).Value = -21200
$ws.Cells.Item(91, 12).Value = -10300
$ws.Cells.Item(91, 13).Value = -10400
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0
$ws.Cells.Item(94, 4).Value = -60900
$ws.Cells.Item(94, 5).Value = -13000
$ws.Cells.Item(94, 6).Value = -32100
$ws.Cells.Item(94, 7).Value = -11800
$ws.Cells.Item(94, 8).Value = -31400
$ws.Cells.Item(94, 9).Value = -7900
$ws.Cells.Item(94, 10).Value = -7200
$ws.Cells.Item(94, 11).Value = -91600
$ws.Cells.Item(94, 12).Value = -10300
$ws.Cells.Item(94, 13).Value = -10000
$ws.Cells.Item(96, 4).Value = -4600
$ws.Cells.Item(96, 5).Value = -4700
$ws.Cells.Item(96, 6).Value = -4700
$ws.Cells.Item(96, 7).Value = -4700
$ws.Cells.Item(96, 8).Value = -4700
$ws.Cells.Item(96, 9).Value = -4700
$ws.Cells.Item(96, 10).Value = -4700
$ws.Cells.Item(96, 11).Value = -4800
$ws.Cells.Item(96, 12).Value = -4900
$ws.Cells.Item(96, 13).Value = -4900
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(100, 4).Value = -133000
$ws.Cells.Item(100, 5).Value = -37600
$ws.Cells.Item(100, 6).Value = -37200
$ws.Cells.Item(100, 7).Value = -45300
$ws.Cells.Item(100, 8).Value = -109800
$ws.Cells.Item(100, 9).Value = -33500
$ws.Cells.Item(100, 10).Value = -19500
$ws.Cells.Item(100, 11).Value = -65600
$ws.Cells.Item(100, 12).Value = -163800
$ws.Cells.Item(100, 13).Value = -8600
$ws.Cells.Item(101, 4).Value = -1400
$ws.Cells.Item(101, 5).Value = -900
$ws.Cells.Item(101, 6).Value = -2900
$ws.Cells.Item(101, 7).Value = 1700
$ws.Cells.Item(101, 8).Value = 300
$ws.Cells.Item(101, 9).Value = 1200
$ws.Cells.Item(101, 10).Value = 1200
$ws.Cells.Item(101, 11).Value = 500
$ws.Cells.Item(101, 12).Value = -1500
$ws.Cells.Item(101, 13).Value = -900
$ws.Cells.Item(102, 4).Value = 9900
$ws.Cells.Item(102, 5).Value = 47100
$ws.Cells.Item(102, 6).Value = -45700
$ws.Cells.Item(102, 7).Value = -114500
$ws.Cells.Item(102, 8).Value = -13100
$ws.Cells.Item(102, 9).Value = 95100
$ws.Cells.Item(102, 10).Value = 81400
$ws.Cells.Item(102, 11).Value = -160600
$ws.Cells.Item(102, 12).Value = -39900
$ws.Cells.Item(102, 13).Value = 61500
